$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Insert a new "config" worksheet right before the existing "varios"
# worksheet (so the tab order becomes: ... prop_mat, config, varios)
# ------------------------------------------------------------------
$variosSheet = $wb.Worksheets.Item("varios")
$config = $wb.Worksheets.Add($variosSheet)
$config.Name = "config"

# ------------------------------------------------------------------
# Populate the "config" sheet with the unit/scale settings that used
# to live inline in "varios".
# ------------------------------------------------------------------

# Header row
$xlCenter = -4108
$config.Range("A1").Value2 = "variable"
$config.Range("B1").Value2 = "valor"
$header = $config.Range("A1:B1")
$header.HorizontalAlignment = $xlCenter
$header.Font.Bold = $true
$header.Font.Name = "Calibri"

# Scale factors used when drawing the deformed shape / diagrams
$config.Range("A2").Value2 = "esc_def"
$config.Range("B2").Value2 = 50

$config.Range("A3").Value2 = "esc_faxial"
$config.Range("B3").Value2 = 0.3

$config.Range("A4").Value2 = "esc_V"
$config.Range("B4").Value2 = 0.3

$config.Range("A5").Value2 = "esc_M"
$config.Range("B5").Value2 = 0.3

# Units
$config.Range("A6").Value2 = "U_FUER"
$config.Range("B6").Value2 = "kN"

$config.Range("A7").Value2 = "U_LONG"
$config.Range("B7").Value2 = "m"

# ------------------------------------------------------------------
# Selection / active-sheet bookkeeping, mirroring the saved state in
# the workbook after the edit.
# ------------------------------------------------------------------

# "prop_mat" ended up with its cursor back at A1
$propMat = $wb.Worksheets.Item("prop_mat")
$null = $propMat.Select()
$null = $propMat.Range("A1").Select()

# "varios" keeps its previous cursor position, it just isn't the
# active tab anymore
$null = $variosSheet.Range("A9").Select()

# "config" is the new active sheet, cursor parked at F18
$null = $config.Select()
$null = $config.Range("F18").Select()
